{"js": "// Replace the three-digit x one-digit multiplication problems/answers\n// in the document's table cells with the new values from the diff.\n// Every original text is unique in the document, so a direct\n// search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"396\u00d77=2772\", \"750\u00d76=4500\"],\n  [\"302\u00d79=2718\", \"116\u00d79=1044\"],\n  [\"700\u00d78=5600\", \"674\u00d74=2696\"],\n  [\"493\u00d72=986\", \"848\u00d79=7632\"],\n  [\"636\u00d77=4452\", \"703\u00d76=4218\"],\n  [\"198\u00d72=396\", \"177\u00d78=1416\"],\n  [\"137\u00d72=274\", \"343\u00d76=2058\"],\n  [\"388\u00d76=2328\", \"767\u00d77=5369\"],\n  [\"672\u00d72=1344\", \"669\u00d73=2007\"],\n  [\"624\u00d76=3744\", \"271\u00d72=542\"],\n  [\"350\u00d78=2800\", \"455\u00d73=1365\"],\n  [\"526\u00d76=3156\", \"812\u00d77=5684\"],\n  [\"761\u00d77=5327\", \"933\u00d79=8397\"],\n  [\"881\u00d73=2643\", \"262\u00d79=2358\"],\n  [\"173\u00d75=865\", \"300\u00d79=2700\"],\n  [\"663\u00d76=3978\", \"804\u00d79=7236\"],\n  [\"338\u00d76=2028\", \"720\u00d73=2160\"],\n  [\"636\u00d75=3180\", \"542\u00d76=3252\"],\n  [\"937\u00d73=2811\", \"149\u00d76=894\"],\n  [\"288\u00d75=1440\", \"604\u00d78=4832\"],\n  [\"636\u00d78=5088\", \"403\u00d77=2821\"],\n  [\"180\u00d72=360\", \"371\u00d79=3339\"],\n  [\"268\u00d78=2144\", \"984\u00d77=6888\"],\n  [\"428\u00d74=1712\", \"294\u00d76=1764\"],\n  [\"222\u00d75=1110\", \"949\u00d75=4745\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit x one-digit multiplication problems/answers\n# in the document's table cells with the new values from the diff.\n# Every original text is unique in the document, so Find/Replace on the\n# whole-document Range for each pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"396\u00d77=2772\", \"750\u00d76=4500\"),\n    @(\"302\u00d79=2718\", \"116\u00d79=1044\"),\n    @(\"700\u00d78=5600\", \"674\u00d74=2696\"),\n    @(\"493\u00d72=986\", \"848\u00d79=7632\"),\n    @(\"636\u00d77=4452\", \"703\u00d76=4218\"),\n    @(\"198\u00d72=396\", \"177\u00d78=1416\"),\n    @(\"137\u00d72=274\", \"343\u00d76=2058\"),\n    @(\"388\u00d76=2328\", \"767\u00d77=5369\"),\n    @(\"672\u00d72=1344\", \"669\u00d73=2007\"),\n    @(\"624\u00d76=3744\", \"271\u00d72=542\"),\n    @(\"350\u00d78=2800\", \"455\u00d73=1365\"),\n    @(\"526\u00d76=3156\", \"812\u00d77=5684\"),\n    @(\"761\u00d77=5327\", \"933\u00d79=8397\"),\n    @(\"881\u00d73=2643\", \"262\u00d79=2358\"),\n    @(\"173\u00d75=865\", \"300\u00d79=2700\"),\n    @(\"663\u00d76=3978\", \"804\u00d79=7236\"),\n    @(\"338\u00d76=2028\", \"720\u00d73=2160\"),\n    @(\"636\u00d75=3180\", \"542\u00d76=3252\"),\n    @(\"937\u00d73=2811\", \"149\u00d76=894\"),\n    @(\"288\u00d75=1440\", \"604\u00d78=4832\"),\n    @(\"636\u00d78=5088\", \"403\u00d77=2821\"),\n    @(\"180\u00d72=360\", \"371\u00d79=3339\"),\n    @(\"268\u00d78=2144\", \"984\u00d77=6888\"),\n    @(\"428\u00d74=1712\", \"294\u00d76=1764\"),\n    @(\"222\u00d75=1110\", \"949\u00d75=4745\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
